# Update cryptocurrency price/volume data in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "64.759.11"
$ws.Cells.Item(2, 5).Value = "  -0.06%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.547.94"
$ws.Cells.Item(3, 5).Value = "  +2.90%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.18%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "599.76"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +2.78%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "135.06"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.79%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "3.546.41"
$ws.Cells.Item(7, 5).Value = "  +2.52%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.02%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.48%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +0.73%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.35%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +1.92%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "4.150.78"
$ws.Cells.Item(13, 5).Value = "  +3.07%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.36%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.554.05"
$ws.Cells.Item(15, 5).Value = "  +3.92%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "27.00"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = "  +2.11%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.74%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "64.654.24"
$ws.Cells.Item(18, 5).Value = "  +0.06%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.02"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +3.66%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.35"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +4.75%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +1.35%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "386.02"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = "  +0.69%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +5.08%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "3.691.93"
$ws.Cells.Item(24, 5).Value = "  +3.23%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "74.10"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = "  +2.67%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.00%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +9.47%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.62"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = "  +4.84%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = "  +0.20%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +4.38%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.35"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = "  +2.43%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "3.553.59"
$ws.Cells.Item(32, 5).Value = "  +2.77%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +21.97%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "23.96"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  +4.34%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +0.05%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +1.65%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "169.71"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -0.04%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "6.91"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = "  +1.85%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +4.51%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "4.98"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  +6.38%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0802"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = "  +4.56%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +2.57%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "26.96"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +17.24%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "42.57"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = "  +1.18%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.21%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "4.45"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = "  +3.11%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +8.50%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.44%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "6.92"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = "  +5.20%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "2.450.06"
$ws.Cells.Item(50, 5).Value = "  +11.17%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +12.72%  "
